# "Volledige data set met juiste namen"
#
# Rename the "Elien Eigen N" tabs (space-separated) to the underscored
# "Elien_Eigen_N" form, and rename the "Elien_CircN" tabs to the fuller
# "Elien_Circular_N" form, to be consistent with the other sheet-name
# families already in the workbook (Yenthe_Circ1..6, Elien_Arrow_1..6, ...).
#
# Each data sheet carries its own tab name again in cell A2 (column A is the
# "ID" column), so that label is refreshed to match every rename.

$wb = $excel.ActiveWorkbook

$renames = @(
    @{ Old = "Elien Eigen 1"; New = "Elien_Eigen_1" },
    @{ Old = "Elien Eigen 2"; New = "Elien_Eigen_2" },
    @{ Old = "Elien Eigen 3"; New = "Elien_Eigen_3" },
    @{ Old = "Elien Eigen 4"; New = "Elien_Eigen_4" },
    @{ Old = "Elien Eigen 5"; New = "Elien_Eigen_5" },
    @{ Old = "Elien Eigen 6"; New = "Elien_Eigen_6" },
    @{ Old = "Elien_Circ1";   New = "Elien_Circular_1" },
    @{ Old = "Elien_Circ2";   New = "Elien_Circular_2" },
    @{ Old = "Elien_Circ3";   New = "Elien_Circular_3" },
    @{ Old = "Elien_Circ4";   New = "Elien_Circular_4" },
    @{ Old = "Elien_Circ5";   New = "Elien_Circular_5" },
    @{ Old = "Elien_Circ6";   New = "Elien_Circular_6" }
)

foreach ($r in $renames) {
    $ws = $wb.Worksheets.Item($r.Old)
    $ws.Name = $r.New
    $ws.Range("A2").Value = $r.New
}

# Re-create the cell selections left behind on a handful of tabs while the
# author was clicking through the workbook doing the renames, ending on the
# "Elien_Circular_4" tab (so it is the one left active/selected).
function Activate-And-Select($sheetName, $cellRef) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Activate() | Out-Null
    $ws.Range($cellRef).Select() | Out-Null
}

Activate-And-Select "Elien_Eigen_2"    "C20"
Activate-And-Select "Elien_Eigen_4"    "A2"
Activate-And-Select "Elien_Eigen_5"    "A3"
Activate-And-Select "Elien_Eigen_6"    "L13"
Activate-And-Select "Elien_Circular_1" "A2"
Activate-And-Select "Elien_Circular_2" "A2"
Activate-And-Select "Elien_Circular_3" "A2"
Activate-And-Select "Elien_Circular_5" "A2"
Activate-And-Select "Elien_Circular_6" "K23"
Activate-And-Select "Elien_Circular_4" "R9"
